# semana 16 de 2024
# Update Esperado (C), Observado (D) and valor p (E) columns for several
# rows in the poisson.xlsx sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> hashtable of column letter -> new value
$changes = @{
    3  = @{ D = 4;  E = 0.02 }
    5  = @{ C = 7;  D = 11;  E = 0.05 }
    7  = @{ C = 1;  D = 80 }
    8  = @{ D = 2;  E = 0.18 }
    9  = @{ D = 0;  E = 1 }
    11 = @{ D = 0;  E = 1 }
    12 = @{ C = 38; D = 34 }
    13 = @{ C = 1;  E = 0.37 }
    14 = @{ C = 0;  D = 5;   E = 0 }
    15 = @{ C = 4;  D = 3;   E = 0.2 }
    16 = @{ C = 13 }
    17 = @{ D = 0;  E = 1 }
    18 = @{ C = 1;  E = 0.37 }
    19 = @{ C = 0;  E = 1 }
    20 = @{ C = 9;  D = 8;   E = 0.13 }
    21 = @{ D = 1 }
    22 = @{ C = 7;  D = 4;   E = 0.09 }
    23 = @{ D = 0;  E = 1 }
    24 = @{ D = 0 }
    25 = @{ D = 0 }
    26 = @{ D = 3 }
    28 = @{ D = 7;  E = 0.1 }
    29 = @{ C = 0;  E = 1 }
    31 = @{ D = 0 }
    32 = @{ D = 1;  E = 0 }
    33 = @{ C = 1;  D = 2;   E = 0.18 }
    34 = @{ C = 8;  D = 4;   E = 0.06 }
    36 = @{ C = 7;  D = 9 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $address = "$col$row"
        $ws.Range($address).Value = $cols[$col]
    }
}
